$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the rotation
$cols = @("D","L","M","N","O","P","Q","R","S","T")

# Capture the current (pre-edit) values for each affected row/column
$rows = @(2,4,5,6,7,8)
$old = @{}
foreach ($r in $rows) {
    $old[$r] = @{}
    foreach ($c in $cols) {
        $old[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Cycle: row 2 <- row 6 <- row 4 <- row 5 <- row 8 <- row 7 <- row 2
$mapping = @{
    2 = 6
    6 = 4
    4 = 5
    5 = 8
    8 = 7
    7 = 2
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value2 = $old[$source][$c]
    }
}
